$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.209.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.488.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.90"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.74"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.30%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.483"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.67"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.386"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.086.68"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.28%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000178"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.491.48"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.178.94"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.41"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.97"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.73"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.43"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.62"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.576"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.628.83"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.38"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.73"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000114"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.43%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.47"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.83%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.515.62"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.27%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.07"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.34"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.83"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.52"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.75"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0779"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.804"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.36"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.17"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.11"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.63"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0256"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.27%  "

# Rows 48-50: reorder coins (Maker, SuiNetwork, Cosmos) -> (SuiNetwork, Cosmos, Maker)
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "SuiNetwork"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.918"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.27%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.75"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.53%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.390.35"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.41%  "
